{"js": "// The commit removes the <w:contextualSpacing .../> element from every\n// paragraph's paragraph-properties (w:pPr) throughout the document body.\n// There is no dedicated \"contextualSpacing\" property exposed on the\n// Word.Paragraph / Word.ParagraphFormat Office.js objects, so we round-trip\n// the body through OOXML: read the body's OOXML, strip every\n// <w:contextualSpacing .../> element (self-closing, any attributes) from\n// it, then write the edited OOXML back in place (Replace) over the whole\n// body. This only removes that one element everywhere; all other markup,\n// relationship ids, hyperlinks, and comments are preserved unchanged.\n\nconst body = context.document.body;\nconst ooxmlResult = body.getOoxml();\nawait context.sync();\n\nlet xml = ooxmlResult.value;\n\n// Remove every <w:contextualSpacing .../> self-closing element regardless\n// of its attributes (the document only ever used w:val=\"0\", but match\n// generically to be robust).\nconst updatedXml = xml.replace(/<w:contextualSpacing\\b[^>]*\\/>/g, \"\");\n\nif (updatedXml !== xml) {\n  body.insertOoxml(updatedXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The commit removes the <w:contextualSpacing .../> element from every\n# paragraph's paragraph-properties (w:pPr) throughout the document body.\n# Word's COM object model has no dedicated ContextualSpacing property that\n# is wired up in this host, so we round-trip the document through OOXML:\n# read the whole document's WordOpenXML, strip every\n# <w:contextualSpacing .../> element (self-closing, any attributes) from\n# it, then write the edited OOXML back in place over the whole document\n# range via InsertXML. This only removes that one element everywhere; all\n# other markup, relationship ids, hyperlinks, and comments are preserved\n# unchanged.\n\n$d = $word.ActiveDocument\n$full = $d.Content\n$xml = $full.WordOpenXML\n\n$updated = [System.Text.RegularExpressions.Regex]::Replace($xml, '<w:contextualSpacing\\b[^>]*/>', '')\n\nif ($updated -ne $xml) {\n    $full.InsertXML($updated)\n}\n"}
